# Update "Weekly Pending Total(Rp)" (C) and "Repayment" (D) values for rows 2-6.
# The "Recovery rate" column E contains formulas (D/C) and will recalculate
# automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1857614943
$ws.Range("D2").Value = 518319157

$ws.Range("C3").Value = 7565664080
$ws.Range("D3").Value = 1901981555

$ws.Range("C4").Value = 4413229239
$ws.Range("D4").Value = 1034612279

$ws.Range("C5").Value = 6613099117
$ws.Range("D5").Value = 1458080341

$ws.Range("C6").Value = 7472398163
$ws.Range("D6").Value = 1638915666

$excel.CalculateFullRebuild()
